# Apply cryptocurrency price/volume updates per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.969.84"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "1.641.61"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.54"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.526"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.54"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0614"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0870"
$ws.Range("E11").Value = "  -2.23%  "
$ws.Range("D12").Value = "1.875.61"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").Value = "1.644.63"
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.565"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.66"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("D17").Value = "27.952.76"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.97"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "0.0₃0723"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.75"
$ws.Range("E22").Value = "  +8.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.40"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.15"
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.12"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D33").Value = "1.464.31"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("E37").Value = "  +2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.932"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.562"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0168"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.35"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("B44").Value = "mCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.46"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.81"
$ws.Range("E45").Value = "  +6.85%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.22"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.37"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").Value = "1.784.56"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.08"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  +6.55%  "
$ws.Range("E51").Value = "  +2.24%  "
